$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) to reflect new ranking/order ---
$ws.Range("A1").Value = 'Datos actualizados a 6 de Mayo de 2020 a las 16:03'
$ws.Range("A74").Value = 'Azerbaiyan'
$ws.Range("A75").Value = 'Croacia'
$ws.Range("A115").Value = 'Kenia'
$ws.Range("A116").Value = 'Maldivas'
$ws.Range("A127").Value = 'Tayikistan'
$ws.Range("A128").Value = 'Estado de Palestina'
$ws.Range("A129").Value = 'Venezuela'
$ws.Range("A130").Value = 'Mauricio'
$ws.Range("A131").Value = 'Isla de Man'
$ws.Range("A132").Value = 'Montenegro'
$ws.Range("A137").Value = 'Cabo Verde'
$ws.Range("A138").Value = 'Islas Feroe'
$ws.Range("A205").Value = 'Seychelles'
$ws.Range("A206").Value = 'Montserrat'

# --- Update statistic values (columns B-H) with the latest report figures ---
$ws.Range("B4").Value = 1238801
$ws.Range("C4").Value = 1168
$ws.Range("D4").Value = 201152
$ws.Range("E4").Value = 965315
$ws.Range("G4").Value = 63
$ws.Range("H4").Value = 72334
$ws.Range("B9").Value = 167239
$ws.Range("C9").Value = 232
$ws.Range("E9").Value = 22846
$ws.Range("B44").Value = 9791
$ws.Range("C44").Value = 114
$ws.Range("D44").Value = 1971
$ws.Range("E44").Value = 7617
$ws.Range("F44").Value = 48
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 203
$ws.Range("D57").Value = 1524
$ws.Range("E57").Value = 3232
$ws.Range("D73").Value = 1577
$ws.Range("E73").Value = 644
$ws.Range("B74").Value = 2127
$ws.Range("C74").Value = 67
$ws.Range("D74").Value = 1536
$ws.Range("E74").Value = 563
$ws.Range("F74").Value = 18
$ws.Range("H74").Value = 28
$ws.Range("B75").Value = 2119
$ws.Range("C75").Value = 7
$ws.Range("D75").Value = 1601
$ws.Range("E75").Value = 433
$ws.Range("F75").Value = 14
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 85
$ws.Range("D79").Value = 1750
$ws.Range("E79").Value = 39
$ws.Range("B114").Value = 608
$ws.Range("C114").Value = 19
$ws.Range("D114").Value = 97
$ws.Range("E114").Value = 470
$ws.Range("F114").Value = 4
$ws.Range("B115").Value = 582
$ws.Range("C115").Value = 47
$ws.Range("D115").Value = 190
$ws.Range("E115").Value = 366
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 26
$ws.Range("B116").Value = 574
$ws.Range("C116").Value = 1
$ws.Range("D116").Value = 20
$ws.Range("E116").Value = 552
$ws.Range("H116").Value = 2
$ws.Range("B127").Value = 379
$ws.Range("C127").Value = 86
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 371
$ws.Range("G127").Value = 3
$ws.Range("H127").Value = 8
$ws.Range("B128").Value = 371
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 127
$ws.Range("E128").Value = 242
$ws.Range("F128").Value = 0
$ws.Range("H128").Value = 2
$ws.Range("B129").Value = 367
$ws.Range("C129").Value = 6
$ws.Range("D129").Value = 164
$ws.Range("E129").Value = 193
$ws.Range("F129").Value = 1
$ws.Range("B130").Value = 332
$ws.Range("D130").Value = 320
$ws.Range("E130").Value = 2
$ws.Range("F130").Value = 3
$ws.Range("H130").Value = 10
$ws.Range("B131").Value = 326
$ws.Range("D131").Value = 271
$ws.Range("E131").Value = 32
$ws.Range("F131").Value = 19
$ws.Range("H131").Value = 23
$ws.Range("B132").Value = 324
$ws.Range("D132").Value = 261
$ws.Range("E132").Value = 55
$ws.Range("F132").Value = 2
$ws.Range("H132").Value = 8
$ws.Range("B137").Value = 191
$ws.Range("C137").Value = 5
$ws.Range("D137").Value = 38
$ws.Range("E137").Value = 151
$ws.Range("H137").Value = 2
$ws.Range("B138").Value = 187
$ws.Range("D138").Value = 185
$ws.Range("E138").Value = 2
$ws.Range("H138").Value = 0
$ws.Range("B172").Value = 58
$ws.Range("C172").Value = 6
$ws.Range("E172").Value = 58
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
